# Apply the edits described by the diff:
#  1. Rename the "Include from SNOMED CT" worksheet to "Include #0".
#  2. Update the "Date" metadata value on the Metadata sheet.
#  3. Insert a new "Jurisdiction" metadata row (with an empty value) right
#     before the "Description" row, pushing Description/Purpose/Copyright/
#     Immutable down by one row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Metadata"
$ws2 = $wb.Worksheets.Item(2)   # "Include from SNOMED CT"

# 1. Rename the second sheet.
$ws2.Name = "Include #0"

# 2. Update the Date property value.
$ws1.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# 3. Insert a new blank row above row 11 ("Description"), shifting the
#    remaining metadata rows down.
$ws1.Range("A11:B11").Insert()

# Match the formatting of the surrounding metadata rows (copy the style
# that row 12 now has - the same style used by every other data row).
$ws1.Range("A12:B12").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row's content.
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""
